$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(565, 1).Value = "outputs/2024-04-28/09-17-50"
$ws.Cells.Item(565, 2).Value = $true
$ws.Cells.Item(565, 3).Value = "mimiciii"
$ws.Cells.Item(565, 4).Value = "mlm"
$ws.Cells.Item(565, 5).Value = "NV"
$ws.Cells.Item(565, 6).Value = "None"
$ws.Cells.Item(565, 7).Value = "descemb_bert"
$ws.Cells.Item(565, 8).Value = "'True"
$ws.Cells.Item(565, 9).Value = "'False"
$ws.Cells.Item(565, 10).Value = "'False"
$ws.Cells.Item(565, 11).Value = 1000
$ws.Cells.Item(565, 12).Value = 0.759
$ws.Cells.Item(565, 13).Value = 0
$ws.Cells.Item(565, 14).Value = 0

$ws.Cells.Item(566, 1).Value = "outputs/2024-04-28/10-31-31"
$ws.Cells.Item(566, 2).Value = $true
$ws.Cells.Item(566, 3).Value = "mimiciii"
$ws.Cells.Item(566, 4).Value = "readmission"
$ws.Cells.Item(566, 5).Value = "VA"
$ws.Cells.Item(566, 6).Value = "descemb_bert"
$ws.Cells.Item(566, 7).Value = "ehr_model"
$ws.Cells.Item(566, 8).Value = "'False"
$ws.Cells.Item(566, 9).Value = "'True"
$ws.Cells.Item(566, 10).Value = "'False"
$ws.Cells.Item(566, 11).Value = 430
$ws.Cells.Item(566, 12).Value = 0.471
$ws.Cells.Item(566, 13).Value = 0.506
$ws.Cells.Item(566, 14).Value = 0.044

$ws.Cells.Item(567, 1).Value = "outputs/2024-04-28/12-56-51"
$ws.Cells.Item(567, 2).Value = $true
$ws.Cells.Item(567, 3).Value = "mimiciii"
$ws.Cells.Item(567, 4).Value = "readmission"
$ws.Cells.Item(567, 5).Value = "DSVA"
$ws.Cells.Item(567, 6).Value = "descemb_bert"
$ws.Cells.Item(567, 7).Value = "ehr_model"
$ws.Cells.Item(567, 8).Value = "'False"
$ws.Cells.Item(567, 9).Value = "'True"
$ws.Cells.Item(567, 10).Value = "'False"
$ws.Cells.Item(567, 11).Value = 1000
$ws.Cells.Item(567, 12).Value = 0.44
$ws.Cells.Item(567, 13).Value = 0.507
$ws.Cells.Item(567, 14).Value = 0.044

$ws.Cells.Item(568, 1).Value = "outputs/2024-04-28/19-05-24"
$ws.Cells.Item(568, 2).Value = $true
$ws.Cells.Item(568, 3).Value = "mimiciii"
$ws.Cells.Item(568, 4).Value = "readmission"
$ws.Cells.Item(568, 5).Value = "DSVA_DPE"
$ws.Cells.Item(568, 6).Value = "descemb_bert"
$ws.Cells.Item(568, 7).Value = "ehr_model"
$ws.Cells.Item(568, 8).Value = "'False"
$ws.Cells.Item(568, 9).Value = "'True"
$ws.Cells.Item(568, 10).Value = "'False"
$ws.Cells.Item(568, 11).Value = 403
$ws.Cells.Item(568, 12).Value = 0.435
$ws.Cells.Item(568, 13).Value = 0.503
$ws.Cells.Item(568, 14).Value = 0.043

$ws.Cells.Item(569, 1).Value = "outputs/2024-04-28/21-32-27"
$ws.Cells.Item(569, 2).Value = $true
$ws.Cells.Item(569, 3).Value = "mimiciii"
$ws.Cells.Item(569, 4).Value = "readmission"
$ws.Cells.Item(569, 5).Value = "VC"
$ws.Cells.Item(569, 6).Value = "descemb_bert"
$ws.Cells.Item(569, 7).Value = "ehr_model"
$ws.Cells.Item(569, 8).Value = "'False"
$ws.Cells.Item(569, 9).Value = "'True"
$ws.Cells.Item(569, 10).Value = "'False"
$ws.Cells.Item(569, 11).Value = 717
$ws.Cells.Item(569, 12).Value = 0.528
$ws.Cells.Item(569, 13).Value = 0.503
$ws.Cells.Item(569, 14).Value = 0.043

$ws.Cells.Item(570, 1).Value = "outputs/2024-04-29/01-07-36"
$ws.Cells.Item(570, 2).Value = $true
$ws.Cells.Item(570, 3).Value = "mimiciii"
$ws.Cells.Item(570, 4).Value = "mortality"
$ws.Cells.Item(570, 5).Value = "VA"
$ws.Cells.Item(570, 6).Value = "descemb_bert"
$ws.Cells.Item(570, 7).Value = "ehr_model"
$ws.Cells.Item(570, 8).Value = "'False"
$ws.Cells.Item(570, 9).Value = "'True"
$ws.Cells.Item(570, 10).Value = "'False"
$ws.Cells.Item(570, 11).Value = 457
$ws.Cells.Item(570, 12).Value = 0.755
$ws.Cells.Item(570, 13).Value = 0.505
$ws.Cells.Item(570, 14).Value = 0.09

$ws.Cells.Item(571, 1).Value = "outputs/2024-04-29/03-41-36"
$ws.Cells.Item(571, 2).Value = $true
$ws.Cells.Item(571, 3).Value = "mimiciii"
$ws.Cells.Item(571, 4).Value = "mortality"
$ws.Cells.Item(571, 5).Value = "DSVA"
$ws.Cells.Item(571, 6).Value = "descemb_bert"
$ws.Cells.Item(571, 7).Value = "ehr_model"
$ws.Cells.Item(571, 8).Value = "'False"
$ws.Cells.Item(571, 9).Value = "'True"
$ws.Cells.Item(571, 10).Value = "'False"
$ws.Cells.Item(571, 11).Value = 423
$ws.Cells.Item(571, 12).Value = 0.732
$ws.Cells.Item(571, 13).Value = 0.509
$ws.Cells.Item(571, 14).Value = 0.091

$ws.Cells.Item(572, 1).Value = "outputs/2024-04-29/06-17-18"
$ws.Cells.Item(572, 2).Value = $false
$ws.Cells.Item(572, 3).Value = "mimiciii"
$ws.Cells.Item(572, 4).Value = "mortality"
$ws.Cells.Item(572, 5).Value = "DSVA_DPE"
$ws.Cells.Item(572, 6).Value = "descemb_bert"
$ws.Cells.Item(572, 7).Value = "ehr_model"
$ws.Cells.Item(572, 8).Value = "'False"
$ws.Cells.Item(572, 9).Value = "'True"
$ws.Cells.Item(572, 10).Value = "'False"
$ws.Cells.Item(572, 11).Value = 9
$ws.Cells.Item(572, 12).Value = 1.45
$ws.Cells.Item(572, 13).Value = 0.503
$ws.Cells.Item(572, 14).Value = 0.089
